$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.48
$ws.Range("J2").Value = 2
$ws.Range("Q2").Value = 1.7
$ws.Range("R2").Value = 2.1
$ws.Range("AW2").Value = 7.5

# Row 3 updates
$ws.Range("G3").Value = 2.05
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 3.9
$ws.Range("J3").Value = 2.88
$ws.Range("L3").Value = 4.75
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("Y3").Value = 9.5
$ws.Range("Z3").Value = 17
$ws.Range("AD3").Value = 6.5
$ws.Range("AH3").Value = 8.5
$ws.Range("AI3").Value = 19
$ws.Range("AN3").Value = 3.75
$ws.Range("AO3").Value = 12
$ws.Range("AS3").Value = 251
$ws.Range("BA3").Value = 151
